$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = 85
$ws.Range("G6").Value = 2539.8
$ws.Range("B10").Value = 28931.5
$ws.Range("F68").Value = 52
$ws.Range("G68").Value = 5986.24
$ws.Range("F70").Value = 22
$ws.Range("G70").Value = 2968.9
$ws.Range("F71").Value = 338
$ws.Range("G71").Value = 21530.6
$ws.Range("F74").Value = 142
$ws.Range("G74").Value = 19929.7
$ws.Range("F86").Value = 71
$ws.Range("G86").Value = 8908.370000000001
$ws.Range("B90").Value = 189530.9
$ws.Range("F151").Value = 94
$ws.Range("G151").Value = 8166.72
$ws.Range("F152").Value = 69
$ws.Range("G152").Value = 6092.01
$ws.Range("B156").Value = 33376.52
$ws.Range("F205").Value = 26
$ws.Range("G205").Value = 9805.639999999999
$ws.Range("B216").Value = 44244.26
$ws.Range("B227").Value = 55373
$ws.Range("E227").Value = 163.62
$ws.Range("F227").Value = -94
$ws.Range("G227").Value = -13562.32
$ws.Range("B228").Value = 63520
$ws.Range("E228").Value = 153.4
$ws.Range("F228").Value = 66
$ws.Range("G228").Value = 9522.48
$ws.Range("F249").Value = 140
$ws.Range("G249").Value = 19294.8
$ws.Range("F255").Value = 584
$ws.Range("G255").Value = 100056.72
$ws.Range("B260").Value = 197871.31
$ws.Range("F280").Value = 138
$ws.Range("G280").Value = 23341.32
$ws.Range("F282").Value = 3
$ws.Range("G282").Value = 161.1
$ws.Range("F289").Value = 3
$ws.Range("G289").Value = 19.71
$ws.Range("F293").Value = 43
$ws.Range("G293").Value = 3023.76
$ws.Range("F294").Value = 34
$ws.Range("G294").Value = 2426.24
$ws.Range("F302").Value = 62
$ws.Range("G302").Value = 13075.18
$ws.Range("B304").Value = 183553.65
$ws.Range("F320").Value = 60
$ws.Range("G320").Value = 4119
$ws.Range("F328").Value = 46
$ws.Range("G328").Value = 1711.66
$ws.Range("B330").Value = 29471.72
$ws.Range("F338").Value = 78
$ws.Range("G338").Value = 1848.6
$ws.Range("F345").Value = 62
$ws.Range("G345").Value = 3807.42
$ws.Range("B346").Value = 27003.35
$ws.Range("F434").Value = 18
$ws.Range("G434").Value = 587.52
$ws.Range("B435").Value = 733
$ws.Range("B473").Value = 64830
$ws.Range("E473").Value = 34.9
$ws.Range("F473").Value = 108
$ws.Range("G473").Value = 3545.64
$ws.Range("B474").Value = 60022
$ws.Range("E474").Value = 37.22
$ws.Range("F474").Value = -113
$ws.Range("G474").Value = -3709.79
$ws.Range("F485").Value = 20
$ws.Range("G485").Value = 3509.4
$ws.Range("B488").Value = 31447.4
$ws.Range("F554").Value = 1
$ws.Range("G554").Value = 37.28
$ws.Range("F555").Value = 26
$ws.Range("G555").Value = 1808.56
$ws.Range("B560").Value = 5551.58
$ws.Range("F565").Value = 0
$ws.Range("G565").Value = 0
$ws.Range("F566").Value = 0
$ws.Range("G566").Value = 0
$ws.Range("F569").Value = 0
$ws.Range("G569").Value = 0
$ws.Range("F570").Value = 0
$ws.Range("G570").Value = 0
$ws.Range("F571").Value = 0
$ws.Range("G571").Value = 0
$ws.Range("B572").Value = 65079
$ws.Range("F572").Value = 6
$ws.Range("G572").Value = 245.22
$ws.Range("B573").Value = 65362
$ws.Range("F573").Value = 20
$ws.Range("G573").Value = 817.4
$ws.Range("F574").Value = 0
$ws.Range("G574").Value = 0
$ws.Range("F576").Value = 0
$ws.Range("G576").Value = 0
$ws.Range("F577").Value = 70
$ws.Range("G577").Value = 3009.3
$ws.Range("F580").Value = 58
$ws.Range("G580").Value = 3305.42
$ws.Range("F581").Value = 0
$ws.Range("G581").Value = 0
$ws.Range("F582").Value = 42
$ws.Range("G582").Value = 2393.58
$ws.Range("B583").Value = 18778.68
$ws.Range("F599").Value = 1804
$ws.Range("G599").Value = 294250.44
$ws.Range("F602").Value = 335
$ws.Range("G602").Value = 48457.75
$ws.Range("B606").Value = 466321.82
$ws.Range("F613").Value = 141
$ws.Range("G613").Value = 22441.56
$ws.Range("B618").Value = 44579.26
$ws.Range("B619").Value = 1845733.6
$ws.Range("B620").Value = 1845733.6
